$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that already have a specialty recorded in column G ("Type of
# Specialty"); everything else gets "none" filled in.
$skipRows = @(11, 15, 17, 18, 19, 24)

for ($r = 2; $r -le 41; $r++) {
    if ($skipRows -contains $r) { continue }
    $ws.Cells.Item($r, 7).Value = "none"
}

# Scroll/select to match the saved view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$ws.Range("G25:G41").Select()
